$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 371.2857
$ws.Range("I103").Value = 350
$ws.Range("J103").Value = 424.5
$ws.Range("K103").Value = 1050
$ws.Range("L103").Value = 1273.5
$ws.Range("M103").Value = -464
$ws.Range("N103").Value = -2445.5

$ws.Range("H135").Value = 699.38464
$ws.Range("I135").Value = 585.4706
$ws.Range("J135").Value = 1474
$ws.Range("K135").Value = 5269.2354
$ws.Range("L135").Value = 13266
$ws.Range("M135").Value = -2734.2354
$ws.Range("N135").Value = -18336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 925.375
$ws.Range("I2").Value = 910.4
$ws.Range("J2").Value = 950.3333
$ws.Range("K2").Value = 910.4
$ws.Range("L2").Value = 950.3333
$ws.Range("M2").Value = -797.4
$ws.Range("N2").Value = -1176.3333

$ws.Range("H32").Value = 5328
$ws.Range("I32").Value = 3017.4312
$ws.Range("K32").Value = 3017.4312
$ws.Range("M32").Value = -2730.4312

$ws.Range("H63").Value = 1218.7142
$ws.Range("I63").Value = 698.5
$ws.Range("J63").Value = 1426.8
$ws.Range("K63").Value = 698.5
$ws.Range("L63").Value = 1426.8
$ws.Range("M63").Value = -12.5
$ws.Range("N63").Value = -2798.8

$ws.Range("H66").Value = 1218.7142
$ws.Range("I66").Value = 698.5
$ws.Range("J66").Value = 1426.8
$ws.Range("K66").Value = 3492.5
$ws.Range("L66").Value = 7134
$ws.Range("M66").Value = -60.5
$ws.Range("N66").Value = -13998

$ws.Range("H74").Value = 49119.04
$ws.Range("I74").Value = 68624.47
$ws.Range("J74").Value = 19860.9
$ws.Range("K74").Value = 68624.47
$ws.Range("L74").Value = 19860.9
$ws.Range("M74").Value = -67750.47
$ws.Range("N74").Value = -21608.9

$ws.Range("H77").Value = 49119.04
$ws.Range("I77").Value = 68624.47
$ws.Range("J77").Value = 19860.9
$ws.Range("K77").Value = 343122.35
$ws.Range("L77").Value = 99304.5
$ws.Range("M77").Value = -338754.35
$ws.Range("N77").Value = -108040.5

$ws.Range("H110").Value = 5730.7
$ws.Range("I110").Value = 6365.4
$ws.Range("J110").Value = 4461.3
$ws.Range("K110").Value = 6365.4
$ws.Range("L110").Value = 4461.3
$ws.Range("M110").Value = -4320.4
$ws.Range("N110").Value = -8551.299999999999

$ws.Range("H116").Value = 925.375
$ws.Range("I116").Value = 910.4
$ws.Range("J116").Value = 950.3333
$ws.Range("K116").Value = 910.4
$ws.Range("L116").Value = 950.3333
$ws.Range("M116").Value = 1383.6
$ws.Range("N116").Value = -5538.3333

$ws.Range("H132").Value = 2094.7273
$ws.Range("I132").Value = 1875.5294
$ws.Range("J132").Value = 3750.889
$ws.Range("K132").Value = 5626.5882
$ws.Range("L132").Value = 11252.667
$ws.Range("M132").Value = -3096.5882
$ws.Range("N132").Value = -16312.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 925.375
$ws.Range("I3").Value = 910.4
$ws.Range("J3").Value = 950.3333
$ws.Range("K3").Value = 910.4
$ws.Range("L3").Value = 950.3333
$ws.Range("M3").Value = -796.4
$ws.Range("N3").Value = -1178.3333

$ws.Range("H94").Value = 1062.5927
$ws.Range("I94").Value = 772.4
$ws.Range("J94").Value = 1891.7142
$ws.Range("K94").Value = 772.4
$ws.Range("L94").Value = 1891.7142
$ws.Range("M94").Value = -321.4
$ws.Range("N94").Value = -2793.7142

$ws.Range("H107").Value = 13843.5
$ws.Range("I107").Value = 1337.25
$ws.Range("K107").Value = 1337.25
$ws.Range("M107").Value = 582.75

$ws.Range("H134").Value = 1469.5205
$ws.Range("I134").Value = 1445.2794
$ws.Range("K134").Value = 4335.8382
$ws.Range("M134").Value = -1800.8382

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 47928.87
$ws.Range("I31").Value = 114288.336
$ws.Range("K31").Value = 114288.336
$ws.Range("M31").Value = -113993.336

$ws.Range("H34").Value = 47928.87
$ws.Range("I34").Value = 114288.336
$ws.Range("K34").Value = 114288.336
$ws.Range("M34").Value = -114086.336

$ws.Range("H58").Value = 191147.61
$ws.Range("I58").Value = 202238.72
$ws.Range("J58").Value = 6295.6665
$ws.Range("K58").Value = 202238.72
$ws.Range("L58").Value = 6295.6665
$ws.Range("M58").Value = -202035.72
$ws.Range("N58").Value = -6701.6665

$ws.Range("H63").Value = 80000
$ws.Range("J63").Value = 80000
$ws.Range("L63").Value = 80000
$ws.Range("N63").Value = -81372

$ws.Range("H66").Value = 80000
$ws.Range("J66").Value = 80000
$ws.Range("L66").Value = 240000
$ws.Range("N66").Value = -246864

$ws.Range("H86").Value = 12723.75
$ws.Range("I86").Value = 12998.5
$ws.Range("J86").Value = 12449
$ws.Range("K86").Value = 12998.5
$ws.Range("L86").Value = 12449
$ws.Range("M86").Value = -11875.5
$ws.Range("N86").Value = -14695

$ws.Range("H89").Value = 12723.75
$ws.Range("I89").Value = 12998.5
$ws.Range("J89").Value = 12449
$ws.Range("K89").Value = 64992.5
$ws.Range("L89").Value = 62245
$ws.Range("M89").Value = -59376.5
$ws.Range("N89").Value = -73477

$ws.Range("H99").Value = 3430.5557
$ws.Range("I99").Value = 3031.6667
$ws.Range("J99").Value = 4228.3335
$ws.Range("K99").Value = 3031.6667
$ws.Range("L99").Value = 4228.3335
$ws.Range("M99").Value = -1533.6667
$ws.Range("N99").Value = -7224.3335

$ws.Range("H103").Value = 3808.5715
$ws.Range("I103").Value = 3808.5715
$ws.Range("K103").Value = 3808.5715
$ws.Range("M103").Value = -2636.5715

$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524

$ws.Range("H126").Value = 3430.5557
$ws.Range("I126").Value = 3031.6667
$ws.Range("J126").Value = 4228.3335
$ws.Range("K126").Value = 9095.000100000001
$ws.Range("L126").Value = 12685.0005
$ws.Range("M126").Value = -6625.000100000001
$ws.Range("N126").Value = -17625.0005

$ws.Range("H132").Value = 3319.647
$ws.Range("I132").Value = 3213.75
$ws.Range("K132").Value = 9641.25
$ws.Range("M132").Value = -7111.25

$ws.Range("H134").Value = 43811
$ws.Range("I134").Value = 13574.111
$ws.Range("K134").Value = 40722.333
$ws.Range("M134").Value = -38187.333

$ws.Range("H136").Value = 191147.61
$ws.Range("I136").Value = 202238.72
$ws.Range("J136").Value = 6295.6665
$ws.Range("K136").Value = 606716.16
$ws.Range("L136").Value = 18886.9995
$ws.Range("M136").Value = -604166.16
$ws.Range("N136").Value = -23986.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 39860.848
$ws.Range("I131").Value = 77810.53999999999
$ws.Range("J131").Value = 1911.1538
$ws.Range("K131").Value = 233431.62
$ws.Range("L131").Value = 5733.4614
$ws.Range("M131").Value = -228391.62
$ws.Range("N131").Value = -15813.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2179.4
$ws.Range("I113").Value = 2223.875
$ws.Range("J113").Value = 2001.5
$ws.Range("K113").Value = 2223.875
$ws.Range("L113").Value = 2001.5
$ws.Range("M113").Value = -53.875
$ws.Range("N113").Value = -6341.5

$ws.Range("H126").Value = 14192.315
$ws.Range("I126").Value = 15450.235
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 46350.705
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -43880.705
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2987.608
$ws.Range("I132").Value = 2447.125
$ws.Range("J132").Value = 4953
$ws.Range("K132").Value = 7341.375
$ws.Range("L132").Value = 14859
$ws.Range("M132").Value = -4811.375
$ws.Range("N132").Value = -19919

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1667.5084
$ws.Range("I136").Value = 1501.6459
$ws.Range("J136").Value = 2391.2727
$ws.Range("K136").Value = 4504.9377
$ws.Range("L136").Value = 7173.8181
$ws.Range("M136").Value = -1954.9377
$ws.Range("N136").Value = -12273.8181
